$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the sensor-reading columns (C:H) for rows 2..21 down by one row, so
# row r+1 receives what used to be in row r. Work bottom-up so we never
# overwrite a value before it has been read.
for ($r = 21; $r -ge 2; $r--) {
    for ($c = 3; $c -le 8; $c++) {
        $v = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($r + 1, $c).Value = $v
    }
}

# New sensor readings for the newly "opened up" row 2 (timestamp/label stay
# as they were - only the measurement columns are new).
$ws.Cells.Item(2, 3).Value = -2.385556173324585
$ws.Cells.Item(2, 4).Value = 2.812312602996826
$ws.Cells.Item(2, 5).Value = -0.1285117015242577
$ws.Cells.Item(2, 6).Value = -0.0219911485910415
$ws.Cells.Item(2, 7).Value = -0.0158824957907199
$ws.Cells.Item(2, 8).Value = 0.0395535230636596

# Drop the final data row (now a duplicate of row 21) to keep 21 data rows.
$ws.Rows.Item(22).Delete()
